$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue 'D2' '60.849.47'
Set-TextValue 'E2' '  -3.39%  '
Set-TextValue 'D3' '2.912.71'
Set-TextValue 'E3' '  -4.07%  '
Set-TextValue 'E4' '  -0.04%  '
Set-TextValue 'D5' '587.43'
Set-TextValue 'E5' '  -1.30%  '
Set-TextValue 'D6' '144.62'
Set-TextValue 'E6' '  -6.12%  '
Set-TextValue 'E8' '  -2.40%  '
Set-TextValue 'D9' '2.911.77'
Set-TextValue 'E9' '  -3.98%  '
Set-TextValue 'D10' '6.69'
Set-TextValue 'E10' '  -2.43%  '
Set-TextValue 'E11' '  -5.52%  '
Set-TextValue 'D12' '0.446'
Set-TextValue 'E12' '  -3.87%  '
Set-TextValue 'E13' '  -4.06%  '
Set-TextValue 'D14' '33.34'
Set-TextValue 'E14' '  -6.89%  '
Set-TextValue 'E15' '  +1.51%  '
Set-TextValue 'D16' '3.394.10'
Set-TextValue 'E16' '  -4.11%  '
Set-TextValue 'D17' '60.806.01'
Set-TextValue 'E17' '  -3.43%  '
Set-TextValue 'E18' '  -5.09%  '
Set-TextValue 'D19' '2.913.67'
Set-TextValue 'E19' '  -4.05%  '
Set-TextValue 'D20' '427.95'
Set-TextValue 'E20' '  -5.77%  '
Set-TextValue 'D21' '13.55'
Set-TextValue 'E21' '  -5.20%  '
Set-TextValue 'D22' '0.679'
Set-TextValue 'E22' '  -2.72%  '
Set-TextValue 'E23' '  -5.71%  '
Set-TextValue 'D24' '80.59'
Set-TextValue 'E24' '  -3.08%  '
Set-TextValue 'D25' '2.22'
Set-TextValue 'E25' '  -3.33%  '
Set-TextValue 'D26' '10.69'
Set-TextValue 'E26' '  -4.96%  '
Set-TextValue 'D27' '11.94'
Set-TextValue 'E27' '  -3.96%  '
Set-TextValue 'E28' '  +0.06%  '
Set-TextValue 'E29' '  -0.02%  '
Set-TextValue 'D30' '7.19'
Set-TextValue 'E30' '  -4.22%  '
Set-TextValue 'E31' '  -3.51%  '
Set-TextValue 'E32' '  -4.06%  '
Set-TextValue 'D33' '26.51'
Set-TextValue 'E33' '  -4.03%  '
Set-TextValue 'E34' '  -4.05%  '
Set-TextValue 'E35' '  +1.73%  '
Set-TextValue 'E36' '  -2.90%  '
Set-TextValue 'E37' '  -5.53%  '
Set-TextValue 'D38' '2.99'
Set-TextValue 'E38' '  -6.69%  '
Set-TextValue 'B39' 'OKB'
Set-TextValue 'C39' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D39' '49.51'
Set-TextValue 'E39' '  -1.67%  '
Set-TextValue 'B40' 'Kaspa'
Set-TextValue 'C40' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D40' '0.126'
Set-TextValue 'E40' '  -3.81%  '
Set-TextValue 'E41' '  -5.66%  '
Set-TextValue 'D42' '8.58'
Set-TextValue 'E42' '  -6.02%  '
Set-TextValue 'D43' '0.296'
Set-TextValue 'E43' '  -3.24%  '
Set-TextValue 'D44' '41.52'
Set-TextValue 'E44' '  -4.96%  '
Set-TextValue 'D45' '377.49'
Set-TextValue 'E45' '  -3.46%  '
Set-TextValue 'E46' '  -3.69%  '
Set-TextValue 'D47' '2.684.20'
Set-TextValue 'E47' '  -1.63%  '
Set-TextValue 'D48' '132.40'
Set-TextValue 'E48' '  -0.26%  '
Set-TextValue 'D50' '24.43'
Set-TextValue 'E50' '  -0.95%  '
Set-TextValue 'E51' '  -2.53%  '
